$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.994.16'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.298.30'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('E7').Value = '  -0.98%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +3.03%  '
$ws.Range('E10').Value = '  +8.40%  '
$ws.Range('E11').Value = '  -0.86%  '
$ws.Range('E13').Value = '  +8.55%  '
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '2.655.33'
$ws.Range('E15').Value = '  +0.04%  '
$ws.Range('D16').Value = '2.333.98'
$ws.Range('E16').Value = '  +1.27%  '
$ws.Range('E17').Value = '  -0.87%  '
$ws.Range('D18').Value = '42.900.34'
$ws.Range('E19').Value = '  +8.60%  '
$ws.Range('E20').Value = '  +0.31%  '
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('E24').Value = '  +10.68%  '
$ws.Range('E25').Value = '  +0.44%  '
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('E28').Value = '  +15.39%  '
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('E34').Value = '  +5.16%  '
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('E36').Value = '  +1.03%  '
$ws.Range('E37').Value = '  -0.57%  '
$ws.Range('E38').Value = '  +0.34%  '
$ws.Range('E39').Value = '  +0.25%  '
$ws.Range('E40').Value = '  +2.23%  '
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +3.93%  '
$ws.Range('E43').Value = '  -2.43%  '
$ws.Range('D44').Value = '1.966.53'
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('E45').Value = '  +3.48%  '
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('E47').Value = '  +0.03%  '
$ws.Range('E48').Value = '  +4.51%  '
$ws.Range('D49').Value = '2.521.96'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  +3.72%  '
$ws.Range('E51').Value = '  +1.24%  '

# Cells whose new value would be auto-coerced to a number by Excel;
# force text storage, then restore the original (default) cell style.
$style = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '300.77'
$ws.Range('D5').Style = $style
$style = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.35'
$ws.Range('D6').Style = $style
$style = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.14'
$ws.Range('D10').Style = $style
$style = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.116'
$ws.Range('D12').Style = $style
$style = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.94'
$ws.Range('D14').Style = $style
$style = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.48'
$ws.Range('D19').Style = $style
$style = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.10'
$ws.Range('D21').Style = $style
$style = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.59'
$ws.Range('D23').Style = $style
$style = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('D24').Style = $style
$style = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('D26').Style = $style
$style = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.97'
$ws.Range('D27').Style = $style
$style = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.39'
$ws.Range('D29').Style = $style
$style = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '167.12'
$ws.Range('D30').Style = $style
$style = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.66'
$ws.Range('D34').Style = $style
$style = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.62'
$ws.Range('D35').Style = $style
$style = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.81'
$ws.Range('D39').Style = $style
$style = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.78'
$ws.Range('D40').Style = $style
$style = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.109'
$ws.Range('D41').Style = $style
$style = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0291'
$ws.Range('D42').Style = $style
$style = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.45'
$ws.Range('D47').Style = $style
$style = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.36'
$ws.Range('D48').Style = $style
$style = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '70.60'
$ws.Range('D51').Style = $style
